$d = $word.ActiveDocument

$d.Content.Find.Execute("82÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷4=", 2) | Out-Null
$d.Content.Find.Execute("34÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷2=", 2) | Out-Null
$d.Content.Find.Execute("51÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷8=", 2) | Out-Null
$d.Content.Find.Execute("15÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷6=", 2) | Out-Null
$d.Content.Find.Execute("11÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷3=", 2) | Out-Null
$d.Content.Find.Execute("74÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷3=", 2) | Out-Null
$d.Content.Find.Execute("28÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷2=", 2) | Out-Null
$d.Content.Find.Execute("40÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "85÷4=", 2) | Out-Null
$d.Content.Find.Execute("86÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "50÷7=", 2) | Out-Null
$d.Content.Find.Execute("71÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷9=", 2) | Out-Null
$d.Content.Find.Execute("20÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷7=", 2) | Out-Null
$d.Content.Find.Execute("97÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷6=", 2) | Out-Null
$d.Content.Find.Execute("62÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷2=", 2) | Out-Null
$d.Content.Find.Execute("96÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "39÷6=", 2) | Out-Null
$d.Content.Find.Execute("49÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "96÷8=", 2) | Out-Null
$d.Content.Find.Execute("14÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷7=", 2) | Out-Null
$d.Content.Find.Execute("27÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "79÷4=", 2) | Out-Null
$d.Content.Find.Execute("34÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷8=", 2) | Out-Null
$d.Content.Find.Execute("28÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷8=", 2) | Out-Null
$d.Content.Find.Execute("62÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷3=", 2) | Out-Null
$d.Content.Find.Execute("68÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷7=", 2) | Out-Null
$d.Content.Find.Execute("83÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "13÷9=", 2) | Out-Null
$d.Content.Find.Execute("60÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷7=", 2) | Out-Null
$d.Content.Find.Execute("86÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "46÷2=", 2) | Out-Null
$d.Content.Find.Execute("58÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷4=", 2) | Out-Null
